$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Meal list income/due entries (rows 3,5,6,8,9) ---
$ws.Range("B3").Value = 2638
$ws.Range("D3").Value = 638

$ws.Range("B5").Value = 3142
$ws.Range("D5").Value = 582

$ws.Range("B6").Value = 723
$ws.Range("C6").Value = 277

$ws.Range("D8").Value = 717

$ws.Range("B9").Value = 557
$ws.Range("C9").Value = 543

# --- Day 17 (17 Oct) meal counts: K17:T17 ---
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 2
$ws.Range("O17").Value = 2
$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 0
$ws.Range("T17").Value = 2

# --- Day 18 (18 Oct): bazar cost + meal counts ---
$ws.Range("F18").Value = 75
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 2
$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 2
$ws.Range("R18").Value = 2
$ws.Range("S18").Value = 0
$ws.Range("T18").Value = 2

# --- Day 19 (19 Oct): bazar cost + meal counts ---
$ws.Range("F19").Value = 1100
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 2
$ws.Range("O19").Value = 1
$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = 2
$ws.Range("S19").Value = 1
$ws.Range("T19").Value = 2

# --- Day 20 (20 Oct): bazar cost + meal counts ---
$ws.Range("F20").Value = 150
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 2
$ws.Range("M20").Value = 2
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 2
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 2
$ws.Range("R20").Value = 2
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 2

# --- H30: blank-ish placeholder string changed from " " (1 space) to "    " (4 spaces) ---
$ws.Range("H30").Value = "    "

# --- Update current selection to reflect where the editor left off ---
$ws.Range("G21").Select()
